$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add new shared-string-backed task rows (23-26) -----------------------
# Rows 23-25 reuse the exact formatting of row 22 (the prior last data row):
#   A: s8, B: s13, C: s38 (date), D: s16, E: (no explicit style), F: s34 (date)
$ws.Range("A22:F22").Copy($ws.Range("A23:F23"))
$ws.Range("A22:F22").Copy($ws.Range("A24:F24"))
$ws.Range("A22:F22").Copy($ws.Range("A25:F25"))
# Row 26 only uses columns A-C (no Observacao/Status/Data entrega), so only
# copy that portion of the formatting across.
$ws.Range("A22:C22").Copy($ws.Range("A26:C26"))

# Row 23
$ws.Cells.Item(23,1).Value = "Alterar arestas para origem e destino ficar na ordem quando trazer a lista do caminho"
$ws.Cells.Item(23,2).Value = "Ivens"
$ws.Cells.Item(23,3).Value = 43261
$ws.Cells.Item(23,4).Value = "entregue dentro do prazo"
$ws.Cells.Item(23,5).Value = "Pronto"
$ws.Cells.Item(23,6).Value = 43261
$ws.Rows(23).RowHeight = 30

# Row 24
$ws.Cells.Item(24,1).Value = "Mostrar mensagem usuário de erro ao carregar arquivo"
$ws.Cells.Item(24,2).Value = "Ivens"
$ws.Cells.Item(24,3).Value = 43261
$ws.Cells.Item(24,4).Value = "entregue dentro do prazo"
$ws.Cells.Item(24,5).Value = "Pronto"
$ws.Cells.Item(24,6).Value = 43261
$ws.Rows(24).RowHeight = 30

# Row 25
$ws.Cells.Item(25,1).Value = "Corrigir erro ao carregar mapa"
$ws.Cells.Item(25,2).Value = "Ivens"
$ws.Cells.Item(25,3).Value = 43261
$ws.Cells.Item(25,4).Value = "entregue dentro do prazo"
$ws.Cells.Item(25,5).Value = "Pronto"
$ws.Cells.Item(25,6).Value = 43261

# Row 26 (only Tarefa / Desenvolvedor / Prazo de entrega are filled in)
$ws.Cells.Item(26,1).Value = "Diagrama de casos de uso"
$ws.Cells.Item(26,2).Value = "Douglas"
$ws.Cells.Item(26,3).Value = 43264

# --- View state: keep selection on the new last "Data entrega" cell -------
$ws.Range("F25").Select()
